# The workbook records one price observation per row, sorted (mostly) by
# date. A new weekly observation was inserted as the new row 207, pushing
# every following row (old 207..297) down by one position (new 208..298).
# The sheet's used range therefore grows from A1:R297 to A1:R298.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 207; Excel shifts rows 207..297 down to
# 208..298 and extends the sheet dimension automatically.
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new weekly observation.
$ws.Range("A207").Value = 8
$ws.Range("B207").Value = "Terminal La Palmera de La Serena"
$ws.Range("C207").Value = "Coquimbo"
$ws.Range("D207").Value = 44839
$ws.Range("E207").Value = 4
$ws.Range("F207").Value = 100112021
$ws.Range("G207").Value = "Ají"
$ws.Range("H207").Value = "Inferno"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 520
$ws.Range("K207").Value = 21000
$ws.Range("L207").Value = 22000
$ws.Range("M207").Value = 21500
$ws.Range("N207").Value = "`$/caja 10 kilos"
$ws.Range("O207").Value = "Región de Arica y Parinacota"
$ws.Range("P207").Value = 2150
$ws.Range("Q207").Value = 10
$ws.Range("R207").Value = "Hortaliza"
